$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores values as text in the source data (e.g. "97.300.86",
# "1.00"), not as numbers. Force a text number-format on each price cell that is
# being updated with a new, numeric-looking value so Excel does not silently
# convert it to a floating point Number (which would lose formatting such as
# trailing zeros / thousands separators) before we write the new value.
$priceRows = @(2,3,5,6,7,8,10,11,12,15,16,18,19,20,21,22,23,24,26,27,28,29,30,31,32,34,35,36,37,38,39,40,41,42,43,45,46,47,48,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "97.475.64"
$ws.Range("D3").Value = "3.600.46"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "244.89"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("D6").Value = "1.76"
$ws.Range("E6").Value = "  +17.38%  "
$ws.Range("D7").Value = "653.26"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").Value = "  +6.49%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "1.06"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").Value = "3.598.59"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "44.74"
$ws.Range("E12").Value = "  +4.26%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "4.268.28"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "97.216.09"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "3.600.41"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "18.27"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "0.529"
$ws.Range("E22").Value = "  +8.14%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "3.50"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "520.38"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").Value = "  +4.68%  "
$ws.Range("D26").Value = "7.00"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").Value = "103.98"
$ws.Range("E27").Value = "  +8.49%  "
$ws.Range("D28").Value = "13.28"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "0.179"
$ws.Range("E29").Value = "  +22.37%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.794.09"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").Value = "12.01"
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "0.189"
$ws.Range("E34").Value = "  +5.67%  "
$ws.Range("D35").Value = "0.986"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("D36").Value = "31.97"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "1.69"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "619.69"
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").Value = "8.81"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "0.156"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").Value = "1.93"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").Value = "0.933"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "6.14"
$ws.Range("E45").Value = "  +6.35%  "
$ws.Range("D46").Value = "0.445"
$ws.Range("E46").Value = "  +43.20%  "
$ws.Range("D47").Value = "0.0447"
$ws.Range("E47").Value = "  +6.80%  "
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "23.66"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "8.67"
$ws.Range("E50").Value = "  +5.24%  "
$ws.Range("D51").Value = "3.31"
